$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "updated" (sheet1): bump the Confirmed/Provisional date markers
# ----------------------------------------------------------------------
$wsUpdated = $wb.Worksheets.Item("updated")
$wsUpdated.Activate()
$wsUpdated.Range("A2").Value2 = 43926
$wsUpdated.Range("B2").Value2 = 43929
$wsUpdated.Range("I18").Select()

# ----------------------------------------------------------------------
# Sheet "totals" (sheet2): revise the 43927 row, drop the hard-coded
# formula on G30, and append the new 43929 row.
# ----------------------------------------------------------------------
$wsTotals = $wb.Worksheets.Item("totals")
$wsTotals.Activate()

$wsTotals.Range("C29").Value2 = 5981
$wsTotals.Range("E29").Value2 = 1472
$wsTotals.Range("F29").Value2 = 224
$wsTotals.Range("G29").Value2 = 223
$wsTotals.Range("H29").Value2 = 1568
$wsTotals.Range("I29").Value2 = 299
$wsTotals.Range("J29").Value2 = 48
$wsTotals.Range("K29").Value2 = "https://www.gov.ie/en/press-release/0f1615-statement-from-the-national-public-health-emergency-team-wednesday-8/"

# G30 used to be a formula "=36+G29" -- it is now a hard value
$wsTotals.Range("G30").Value2 = 210

# New row 31 for 2020-04-09 (serial 43929)
$wsTotals.Range("A30").Copy()
$wsTotals.Range("A31").PasteSpecial(-4122)
$wsTotals.Range("A31").Value2 = 43929

$wsTotals.Range("B30").Copy()
$wsTotals.Range("B31").PasteSpecial(-4122)
$wsTotals.Range("B31").Value2 = "Provisional"

$wsTotals.Range("C30").Copy()
$wsTotals.Range("C31").PasteSpecial(-4122)
$wsTotals.Range("C31").Value2 = 6074

$wsTotals.Range("D30").Copy()
$wsTotals.Range("D31").PasteSpecial(-4122)
$wsTotals.Range("D31").Formula = "=C31-C30"

$wsTotals.Range("E31").Value2 = "NA"
$wsTotals.Range("F31").Value2 = "NA"

$wsTotals.Range("G30").Copy()
$wsTotals.Range("G31").PasteSpecial(-4122)
$wsTotals.Range("G31").Value2 = 235

$wsTotals.Range("H31").Value2 = "NA"
$wsTotals.Range("I31").Value2 = "NA"
$wsTotals.Range("J31").Value2 = "NA"
$wsTotals.Range("K31").Value2 = "https://www.gov.ie/en/press-release/0f1615-statement-from-the-national-public-health-emergency-team-wednesday-8/"

$wsTotals.Rows.Item(31).RowHeight = 17

$winTotals = $excel.ActiveWindow
$winTotals.ScrollRow = 2
$winTotals.ScrollColumn = 2
$wsTotals.Range("G30").Select()

Write-Output "done"
